$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: id / itemType / key (name -> key localization change) ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "itemType"

# --- Data rows: id, itemType, key(numeric) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Item_1"
$ws.Range("C2").Value = 10000

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Item_2"
$ws.Range("C3").Value = 10001

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Item_PuzzleBlock_A"
$ws.Range("C4").Value = 10100

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Item_PuzzleBlock_B"
$ws.Range("C5").Value = 10101

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Item_Key_A"
$ws.Range("C6").Value = 10201

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Item_Book_A"
$ws.Range("C7").Value = 10301

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Item_Book_B"
$ws.Range("C8").Value = 10302

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Item_Book_C"
$ws.Range("C9").Value = 10303

# Header C1 written last so the shared-string table order matches
$ws.Range("C1").Value = "key"

# --- Conditional-format-like cell styles ("Good" = green, "Neutral" = yellow) ---
# Item_1, Item_2, Item_Key_A rows -> Good
$ws.Range("A2:C2").Style = "Good"
$ws.Range("A3:C3").Style = "Good"
$ws.Range("A6:C6").Style = "Good"

# PuzzleBlock_A/B, Book_A/B/C rows -> Neutral
$ws.Range("A4:C4").Style = "Neutral"
$ws.Range("A5:C5").Style = "Neutral"
$ws.Range("A7:C7").Style = "Neutral"
$ws.Range("A8:C8").Style = "Neutral"
$ws.Range("A9:C9").Style = "Neutral"

# --- Column widths (best-fit-like, approximated) ---
$ws.Columns("B").ColumnWidth = 16.666666666666668
$ws.Columns("C").ColumnWidth = 11.5
$ws.Columns("D").ColumnWidth = 41.666666666666664
$ws.Columns("E").ColumnWidth = 13.333333333333334
$ws.Columns("F").ColumnWidth = 39.0

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("F12").Select()
